# Updated symbol list on Sun Jan 22 03:34:54 UTC 2023 with GitHub Actions
# Refresh crypto Price (column D) and Volume(1h) (column E) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  D="300.60";       E="-0.53%"},
    @{Row=3;  D="38.15";        E="8.55%"},
    @{Row=4;  D="4.997";        E="-2.94%"},
    @{Row=5;                    E="-0.58%"},
    @{Row=6;  D="2.172";        E="-7.65%"},
    @{Row=7;  D="7.968";        E="-0.85%"},
    @{Row=8;  D="3.993";        E="1.25%"},
    @{Row=9;  D="0.9168";       E="-1.58%"},
    @{Row=10; D="0.09063";      E="-8.98%"},
    @{Row=11; D="0.1793";       E="-0.01%"},
    @{Row=12; D="0.08434";      E="-1.92%"},
    @{Row=13; D="0.03540";      E="6.96%"},
    @{Row=14; D="0.09935";      E="0.16%"},
    @{Row=15; D="0.001480";     E="-1.35%"},
    @{Row=16; D="0.005687";     E="-1.22%"},
    @{Row=17; D="3.476";        E="0.44%"},
    @{Row=18; D="2.223";        E="3.64%"},
    @{Row=20; D="0.1317";       E="1.08%"},
    @{Row=21; D="4.567";        E="6.25%"},
    @{Row=22; D="0.2235";       E="-2.95%"},
    @{Row=23; D="0.04659";      E="2.74%"},
    @{Row=24; D="0.001229";     E="1.19%"},
    @{Row=25; D="0.004438";     E="1.47%"},
    @{Row=26; D="0.0001302";    E="0.03%"},
    @{Row=27; D="0.0004756";    E="40.09%"},
    @{Row=39; D="0.01743";      E="-2.97%"},
    @{Row=40; D="0.04687";      E="-2.26%"},
    @{Row=41; D="0.007907";     E="1.64%"},
    @{Row=42; D="0.1387";       E="-1.72%"},
    @{Row=43; D="0.007688";     E="12.62%"},
    @{Row=44; D="0.002293";     E="10.61%"},
    @{Row=45; D="0.009765";     E="3.42%"},
    @{Row=46; D="0.00006040";   E="-1.27%"},
    @{Row=47; D="0.00000000751";E="0.03%"},
    @{Row=48; D="8.587";        E="187.70%"},
    @{Row=49;                   E="34.80%"},
    @{Row=50; D="0.00002103";   E="0.03%"},
    @{Row=51; D="0.0002003";    E="0.03%"}
)

foreach ($item in $data) {
    $r = $item.Row

    if ($item.ContainsKey("D")) {
        $cellD = $ws.Range("D$r")
        $cellD.NumberFormat = "@"
        $cellD.Value = $item.D
    }

    if ($item.ContainsKey("E")) {
        $cellE = $ws.Range("E$r")
        $cellE.NumberFormat = "@"
        $cellE.Value = $item.E
    }
}
